$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Mean" summary row (row 25) below the data table (rows 3:22).
$ws.Range("B25").Value = "Mean"
$ws.Range("C25").Formula = "=AVERAGE(C3:C22)"
$ws.Range("D25:G25").Formula = "=AVERAGE(D3:D22)"

# Leave the selection where the author last left it when saving.
$ws.Range("F31").Select() | Out-Null
